$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(806).Insert()

$ws.Cells.Item(806, 1).Value = 11
$ws.Cells.Item(806, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(806, 3).Value = "Bíobío"
$ws.Cells.Item(806, 4).Value = Get-Date -Year 2023 -Month 12 -Day 5 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(806, 5).Value = 8
$ws.Cells.Item(806, 6).Value = "Fruta"
$ws.Cells.Item(806, 7).Value = 100104
$ws.Cells.Item(806, 8).Value = "Frutos de pepita"
$ws.Cells.Item(806, 9).Value = 100104005
$ws.Cells.Item(806, 10).Value = "Pera"
$ws.Cells.Item(806, 11).Value = "Packham's Triumph"
$ws.Cells.Item(806, 12).Value = "Primera"
$ws.Cells.Item(806, 13).Value = 220
$ws.Cells.Item(806, 14).Value = 15000
$ws.Cells.Item(806, 15).Value = 16000
$ws.Cells.Item(806, 16).Value = 15455
$ws.Cells.Item(806, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(806, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(806, 19).Value = 859
$ws.Cells.Item(806, 20).Value = 18
